$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 986.35
$ws.Range("I70").Value = 903.7
$ws.Range("J70").Value = 1069
$ws.Range("K70").Value = 2711.1
$ws.Range("L70").Value = 3207
$ws.Range("M70").Value = -2441.1
$ws.Range("N70").Value = -3747

# Row 73
$ws.Range("H73").Value = 986.35
$ws.Range("I73").Value = 903.7
$ws.Range("J73").Value = 1069
$ws.Range("K73").Value = 2711.1
$ws.Range("L73").Value = 3207
$ws.Range("M73").Value = -1775.1
$ws.Range("N73").Value = -5079

# Row 112
$ws.Range("H112").Value = 2206.4443
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2206.4443
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6619.3329
$ws.Range("N112").Value = -8835.332900000001
$ws.Range("M112").ClearContents()

# Row 129
$ws.Range("H129").Value = 883.6087
$ws.Range("I129").Value = 754
$ws.Range("J129").Value = 952.73334
$ws.Range("K129").Value = 2262
$ws.Range("L129").Value = 2858.20002
$ws.Range("M129").Value = 2738
$ws.Range("N129").Value = -12858.20002

# Row 132
$ws.Range("H132").Value = 11120902
$ws.Range("I132").Value = 14499558
$ws.Range("J132").Value = 19600.857
$ws.Range("K132").Value = 43498674
$ws.Range("L132").Value = 58802.571
$ws.Range("M132").Value = -43496144
$ws.Range("N132").Value = -63862.571

# Row 137
$ws.Range("H137").Value = 6288.6665
$ws.Range("I137").Value = 12000
$ws.Range("J137").Value = 3433
$ws.Range("K137").Value = 36000
$ws.Range("L137").Value = 10299
$ws.Range("M137").Value = -33450
$ws.Range("N137").Value = -15399

# Row 138
$ws.Range("H138").Value = 2908.405
$ws.Range("I138").Value = 2799.5715
$ws.Range("J138").Value = 2918.986
$ws.Range("K138").Value = 8398.7145
$ws.Range("L138").Value = 8756.957999999999
$ws.Range("M138").Value = -3258.7145
$ws.Range("N138").Value = -19036.958

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11943.14
$ws.Range("I32").Value = 9201.429
$ws.Range("J32").Value = 16092.757
$ws.Range("K32").Value = 9201.429
$ws.Range("L32").Value = 16092.757
$ws.Range("M32").Value = -8914.429
$ws.Range("N32").Value = -16666.757

# Row 74
$ws.Range("H74").Value = 1799.8148
$ws.Range("I74").Value = 994.94116
$ws.Range("J74").Value = 3168.1
$ws.Range("K74").Value = 994.94116
$ws.Range("L74").Value = 3168.1
$ws.Range("M74").Value = -120.94116
$ws.Range("N74").Value = -4916.1

# Row 77
$ws.Range("H77").Value = 1799.8148
$ws.Range("I77").Value = 994.94116
$ws.Range("J77").Value = 3168.1
$ws.Range("K77").Value = 4974.7058
$ws.Range("L77").Value = 15840.5
$ws.Range("M77").Value = -606.7057999999997
$ws.Range("N77").Value = -24576.5

# Row 110
$ws.Range("H110").Value = 824.2727
$ws.Range("I110").Value = 806.7
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 806.7
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1238.3
$ws.Range("N110").Value = -5090

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 50004336
$ws.Range("I86").Value = 58827836
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 58827836
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -58826713
$ws.Range("N86").Value = -6746

# Row 89
$ws.Range("H89").Value = 50004336
$ws.Range("I89").Value = 58827836
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 294139180
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -294133564
$ws.Range("N89").Value = -33732

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 54113.383
$ws.Range("I22").Value = 277.4
$ws.Range("J22").Value = 233566.67
$ws.Range("K22").Value = 277.4
$ws.Range("L22").Value = 233566.67
$ws.Range("M22").Value = 72.60000000000002
$ws.Range("N22").Value = -234266.67

# Row 31
$ws.Range("H31").Value = 2205.4856
$ws.Range("I31").Value = 2123.5881
$ws.Range("J31").Value = 4990
$ws.Range("K31").Value = 2123.5881
$ws.Range("L31").Value = 4990
$ws.Range("M31").Value = -1828.5881
$ws.Range("N31").Value = -5580

# Row 34
$ws.Range("H34").Value = 2205.4856
$ws.Range("I34").Value = 2123.5881
$ws.Range("J34").Value = 4990
$ws.Range("K34").Value = 2123.5881
$ws.Range("L34").Value = 4990
$ws.Range("M34").Value = -1921.5881
$ws.Range("N34").Value = -5394

# Row 62
$ws.Range("H62").Value = 16668967
$ws.Range("I62").Value = 2480
$ws.Range("J62").Value = 100001400
$ws.Range("K62").Value = 2480
$ws.Range("L62").Value = 100001400
$ws.Range("M62").Value = -1856
$ws.Range("N62").Value = -100002648

# Row 65
$ws.Range("H65").Value = 16668967
$ws.Range("I65").Value = 2480
$ws.Range("J65").Value = 100001400
$ws.Range("K65").Value = 12400
$ws.Range("L65").Value = 500007000
$ws.Range("M65").Value = -9280
$ws.Range("N65").Value = -500013240

# Row 68
$ws.Range("H68").Value = 14750
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 14750
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 14750
$ws.Range("N68").Value = -16248
$ws.Range("M68").ClearContents()

# Row 71
$ws.Range("H71").Value = 14750
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 14750
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 44250
$ws.Range("N71").Value = -51738
$ws.Range("M71").ClearContents()

# Row 74
$ws.Range("H74").Value = 32999.668
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 32999.668
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 32999.668
$ws.Range("N74").Value = -34747.668

# Row 77
$ws.Range("H77").Value = 32999.668
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 32999.668
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 98999.00399999999
$ws.Range("N77").Value = -107735.004

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 9092111
$ws.Range("I34").Value = 228.6
$ws.Range("J34").Value = 16668680
$ws.Range("K34").Value = 685.8
$ws.Range("L34").Value = 50006040
$ws.Range("M34").Value = -601.8
$ws.Range("N34").Value = -50006208

# Row 68
$ws.Range("H68").Value = 936.46155
$ws.Range("I68").Value = 996
$ws.Range("J68").Value = 899.25
$ws.Range("K68").Value = 2988
$ws.Range("L68").Value = 2697.75
$ws.Range("M68").Value = -2177
$ws.Range("N68").Value = -4319.75

# Row 71
$ws.Range("H71").Value = 936.46155
$ws.Range("I71").Value = 996
$ws.Range("J71").Value = 899.25
$ws.Range("K71").Value = 8964
$ws.Range("L71").Value = 8093.25
$ws.Range("M71").Value = -4908
$ws.Range("N71").Value = -16205.25

# Row 131
$ws.Range("H131").Value = 23846266
$ws.Range("I131").Value = 142857630
$ws.Range("J131").Value = 43994.17
$ws.Range("K131").Value = 428572890
$ws.Range("L131").Value = 131982.51
$ws.Range("M131").Value = -428567850
$ws.Range("N131").Value = -142062.51

# Row 133
$ws.Range("H133").Value = 6040.2085
$ws.Range("I133").Value = 1000
$ws.Range("J133").Value = 6259.3477
$ws.Range("K133").Value = 3000
$ws.Range("L133").Value = 18778.0431
$ws.Range("M133").Value = 2060
$ws.Range("N133").Value = -28898.0431

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1736.6666
$ws.Range("I113").Value = 1736.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1736.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 433.3334
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2350.25
$ws.Range("I22").Value = 5001
$ws.Range("J22").Value = 1466.6666
$ws.Range("K22").Value = 5001
$ws.Range("L22").Value = 1466.6666
$ws.Range("M22").Value = -4706
$ws.Range("N22").Value = -2056.6666

# Row 27
$ws.Range("H27").Value = 2350.25
$ws.Range("I27").Value = 5001
$ws.Range("J27").Value = 1466.6666
$ws.Range("K27").Value = 5001
$ws.Range("L27").Value = 1466.6666
$ws.Range("M27").Value = -4894
$ws.Range("N27").Value = -1680.6666

# Row 61
$ws.Range("H61").Value = 1120.7858
$ws.Range("I61").Value = 1120.7858
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1120.7858
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -918.7858000000001
$ws.Range("N61").ClearContents()

# Row 113
$ws.Range("H113").Value = 1120.7858
$ws.Range("I113").Value = 1120.7858
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1120.7858
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1049.2142
$ws.Range("N113").ClearContents()

# Row 136
$ws.Range("H136").Value = 1936.3636
$ws.Range("I136").Value = 1830
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5490
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2940
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1837.25
$ws.Range("I81").Value = 1400.25
$ws.Range("J81").Value = 1910.0834
$ws.Range("K81").Value = 2800.5
$ws.Range("L81").Value = 3820.1668
$ws.Range("M81").Value = -1739.5
$ws.Range("N81").Value = -5942.1668

# Row 84
$ws.Range("H84").Value = 1837.25
$ws.Range("I84").Value = 1400.25
$ws.Range("J84").Value = 1910.0834
$ws.Range("K84").Value = 14002.5
$ws.Range("L84").Value = 19100.834
$ws.Range("M84").Value = -8698.5
$ws.Range("N84").Value = -29708.834

# Row 113
$ws.Range("H113").Value = 479.8
$ws.Range("I113").Value = 228.28572
$ws.Range("J113").Value = 1066.6666
$ws.Range("K113").Value = 684.85716
$ws.Range("L113").Value = 3199.9998
$ws.Range("M113").Value = 1485.14284
$ws.Range("N113").Value = -7539.9998

# Row 136
$ws.Range("H136").Value = 1453.6938
$ws.Range("I136").Value = 557.9655
$ws.Range("J136").Value = 2752.5
$ws.Range("K136").Value = 1673.8965
$ws.Range("L136").Value = 8257.5
$ws.Range("M136").Value = 876.1034999999999
$ws.Range("N136").Value = -13357.5
